$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: status text updated for both locales (zh-cn / de-de columns)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in handback columns now that the handback is done
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Value = "2568e378-3ee1-41c0-928b-7338821fda23.md"
$wsZh.Range("J2").Value = "2568e378-3ee1-41c0-928b-7338821fda23.57846a3a715bc65426d34efc1ef3c16cc29fdb49.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-17 12:27:58"

$wsZh.Range("I3").Value = "bf361fd3-ec84-44d6-9753-1c17e726b73e.md"
$wsZh.Range("J3").Value = "bf361fd3-ec84-44d6-9753-1c17e726b73e.154e316aaa4e5f3ab9ca6558e1321661915e27bf.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-17 12:27:58"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3645068452f80448ad0d9a055fb472fa56c31efb/e2e/2568e378-3ee1-41c0-928b-7338821fda23.md", "", "", "2568e378-3ee1-41c0-928b-7338821fda23.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3645068452f80448ad0d9a055fb472fa56c31efb/e2e/bf361fd3-ec84-44d6-9753-1c17e726b73e.md", "", "", "bf361fd3-ec84-44d6-9753-1c17e726b73e.md")

$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I2").Font.Name = "Calibri"
$wsZh.Range("I2").Font.Size = 11
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = 15570276
$wsZh.Range("I3").Font.Name = "Calibri"
$wsZh.Range("I3").Font.Size = 11

$wsZh.Columns.Item(3).ColumnWidth = 29.09
$wsZh.Columns.Item(9).ColumnWidth = 39.18
$wsZh.Columns.Item(10).ColumnWidth = 39.18

# ---------------------------------------------------------------------------
# de-de sheet: fill in handback columns now that the handback is done
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Value = "2568e378-3ee1-41c0-928b-7338821fda23.md"
$wsDe.Range("J2").Value = "2568e378-3ee1-41c0-928b-7338821fda23.57846a3a715bc65426d34efc1ef3c16cc29fdb49.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-17 12:28:13"

$wsDe.Range("I3").Value = "bf361fd3-ec84-44d6-9753-1c17e726b73e.md"
$wsDe.Range("J3").Value = "bf361fd3-ec84-44d6-9753-1c17e726b73e.154e316aaa4e5f3ab9ca6558e1321661915e27bf.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-17 12:28:13"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3645068452f80448ad0d9a055fb472fa56c31efb/e2e/2568e378-3ee1-41c0-928b-7338821fda23.md", "", "", "2568e378-3ee1-41c0-928b-7338821fda23.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3645068452f80448ad0d9a055fb472fa56c31efb/e2e/bf361fd3-ec84-44d6-9753-1c17e726b73e.md", "", "", "bf361fd3-ec84-44d6-9753-1c17e726b73e.md")

$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I2").Font.Name = "Calibri"
$wsDe.Range("I2").Font.Size = 11
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = 15570276
$wsDe.Range("I3").Font.Name = "Calibri"
$wsDe.Range("I3").Font.Size = 11

$wsDe.Columns.Item(3).ColumnWidth = 29.09
$wsDe.Columns.Item(9).ColumnWidth = 39.18
$wsDe.Columns.Item(10).ColumnWidth = 39.18

# ---------------------------------------------------------------------------
# Overview sheet column widths (E/F grow to fit the longer status text)
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.09
$wsOverview.Columns.Item(6).ColumnWidth = 29.09
